# #5: cash & deposit done
# Extend the "存款" (deposit) sheet with bank / deposit_type / currency
# columns up front and property_category / category / date /
# legislator_name / legislator_id / source_file / index columns at the
# end, matching the layout already used on the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# Style reference cells that already carry the correct look:
#   Cells(1, 2) -> header style (bold, bordered, centered)
#   Cells(2, 2) -> plain data style
$headerStyleSrc = $ws.Cells.Item(1, 2)
$dataStyleSrc = $ws.Cells.Item(2, 2)

# ---- Header row -----------------------------------------------------
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

foreach ($col in 7..13) {
    $headerStyleSrc.Copy()
    $ws.Cells.Item(1, $col).PasteSpecial(-4122)
}

# ---- Data rows --------------------------------------------------------
$rows = @(
    @{ Row = 2; Index = 52; Bank = "第一商業銀行土城分行"; Owner = "吳麗香"; Total = 2446717 },
    @{ Row = 3; Index = 53; Bank = "上城農會土城分行";     Owner = "盧嘉辰"; Total = 1937829 },
    @{ Row = 4; Index = 54; Bank = "中華郵政股份有限公司"; Owner = "盧嘉辰"; Total = 347674 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 2).Value = $r.Bank
    $ws.Cells.Item($row, 3).Value = "活期存款"
    $ws.Cells.Item($row, 4).Value = "新臺幣"
    $ws.Cells.Item($row, 5).Value = $r.Owner
    $ws.Cells.Item($row, 6).Value = $r.Total
    $ws.Cells.Item($row, 7).Value = "deposit"
    $ws.Cells.Item($row, 8).Value = "normal"

    # "2012-04-12" reads as a date literal, so Excel would otherwise
    # silently convert it to a date serial number; force text format
    # first so it round-trips as the original string.
    $dateCell = $ws.Cells.Item($row, 9)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2012-04-12"

    $ws.Cells.Item($row, 10).Value = "盧嘉辰"
    $ws.Cells.Item($row, 11).Value = 1715
    $ws.Cells.Item($row, 12).Value = "tmp79201"
    $ws.Cells.Item($row, 13).Value = $r.Index

    foreach ($col in 7..13) {
        $dataStyleSrc.Copy()
        $ws.Cells.Item($row, $col).PasteSpecial(-4122)
    }
}
